$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32, shifting existing rows 32:44 down to 33:45.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly record.
$ws.Cells.Item(32, 1).Value = 8
$ws.Cells.Item(32, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44736
$ws.Cells.Item(32, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = 100114007
$ws.Cells.Item(32, 7).Value = "Jengibre"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 400
$ws.Cells.Item(32, 11).Value = 16000
$ws.Cells.Item(32, 12).Value = 17000
$ws.Cells.Item(32, 13).Value = 16500
$ws.Cells.Item(32, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(32, 15).Value = "Perú"
$ws.Cells.Item(32, 16).Value = 1269
$ws.Cells.Item(32, 17).Value = 13
$ws.Cells.Item(32, 18).Value = "Hortaliza"
